# Rename existing sheet "Data" -> "production", and add a new "staging"
# sheet (duplicate of production's data, with the Quiz Dashboard Classes
# value for Admin/Principal rows swapped to the new repeated-class string).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "production"

# New sheet placed right after "production"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "staging"

$ws2.Range('A1').Value = 'TCID'
$ws2.Range('B1').Value = 'Platform'
$ws2.Range('C1').Value = 'Method Name'
$ws2.Range('D1').Value = 'Role'
$ws2.Range('E1').Value = 'Key'
$ws2.Range('F1').Value = 'Value'
$ws2.Range('B2').Value = 'Web'
$ws2.Range('C2').Value = 'verifyPrimeClasses_Old'
$ws2.Range('D2').Value = 'Admin'
$ws2.Range('E2').Value = 'Prime Classes'
$ws2.Range('F2').Value = 'Pre Nursery, Nursery, KG, Class 1, Class 2, Class 3, Class 4, Class 5, Class 6, Class 7, Class 8, Class 9, Class 10, Class 11, Class 12'
$ws2.Range('B3').Value = 'Android'
$ws2.Range('C3').Value = 'verifyPrimeClasses_Old'
$ws2.Range('D3').Value = 'Admin'
$ws2.Range('E3').Value = 'Prime Classes'
$ws2.Range('F3').Value = 'Pre Nursery, Nursery, KG, Class 1, Class 2, Class 3, Class 4, Class 5, Class 6, Class 7, Class 8, Class 9, Class 10'
$ws2.Range('B4').Value = 'iOS'
$ws2.Range('C4').Value = 'verifyPrimeClasses_Old'
$ws2.Range('D4').Value = 'Admin'
$ws2.Range('E4').Value = 'Prime Classes'
$ws2.Range('F4').Value = 'Pre Nursery, Nursery, KG, Class 1, Class 2, Class 3, Class 4, Class 5, Class 6, Class 7, Class 8, Class 9, Class 10, Class 11, Class 12'
$ws2.Range('B5').Value = 'Web'
$ws2.Range('C5').Value = 'verifyPrimeClasses_Old'
$ws2.Range('D5').Value = 'Principal'
$ws2.Range('E5').Value = 'Prime Classes'
$ws2.Range('F5').Value = 'Pre Nursery, Nursery, KG, Class 1, Class 2, Class 3, Class 4, Class 5, Class 6, Class 7, Class 8, Class 9, Class 10, Class 11, Class 12'
$ws2.Range('B6').Value = 'Android'
$ws2.Range('C6').Value = 'verifyPrimeClasses_Old'
$ws2.Range('D6').Value = 'Principal'
$ws2.Range('E6').Value = 'Prime Classes'
$ws2.Range('F6').Value = 'Pre Nursery, Nursery, KG, Class 1, Class 2, Class 3, Class 4, Class 5, Class 6, Class 7, Class 8, Class 9, Class 10'
$ws2.Range('B7').Value = 'iOS'
$ws2.Range('C7').Value = 'verifyPrimeClasses_Old'
$ws2.Range('D7').Value = 'Principal'
$ws2.Range('E7').Value = 'Prime Classes'
$ws2.Range('F7').Value = 'Pre Nursery, Nursery, KG, Class 1, Class 2, Class 3, Class 4, Class 5, Class 6, Class 7, Class 8, Class 9, Class 10, Class 11, Class 12'
$ws2.Range('B8').Value = 'Web'
$ws2.Range('C8').Value = 'verifyPrimeClasses_Old'
$ws2.Range('D8').Value = 'Teacher'
$ws2.Range('E8').Value = 'Prime Classes'
$ws2.Range('F8').Value = 'Pre Nursery, Nursery, KG, Class 1, Class 2, Class 3, Class 4, Class 5, Class 6, Class 7, Class 8, Class 9, Class 10, Class 11, Class 12'
$ws2.Range('B9').Value = 'Android'
$ws2.Range('C9').Value = 'verifyPrimeClasses_Old'
$ws2.Range('D9').Value = 'Teacher'
$ws2.Range('E9').Value = 'Prime Classes'
$ws2.Range('F9').Value = 'Pre Nursery, Nursery, KG, Class 1, Class 2, Class 3, Class 4, Class 5, Class 6, Class 7, Class 8, Class 9, Class 10'
$ws2.Range('B10').Value = 'iOS'
$ws2.Range('C10').Value = 'verifyPrimeClasses_Old'
$ws2.Range('D10').Value = 'Teacher'
$ws2.Range('E10').Value = 'Prime Classes'
$ws2.Range('F10').Value = 'Pre Nursery, Nursery, KG, Class 1, Class 2, Class 3, Class 4, Class 5, Class 6, Class 7, Class 8, Class 9, Class 10, Class 11, Class 12'
$ws2.Range('B11').Value = 'Web'
$ws2.Range('C11').Value = 'verifyPrimeClasses_Parent'
$ws2.Range('D11').Value = 'Parent'
$ws2.Range('E11').Value = 'Prime Classes'
$ws2.Range('F11').Value = 'Class 5, Class 6, Class 7'
$ws2.Range('B12').Value = 'Android'
$ws2.Range('C12').Value = 'verifyPrimeClasses_Parent'
$ws2.Range('D12').Value = 'Parent'
$ws2.Range('E12').Value = 'Prime Classes'
$ws2.Range('F12').Value = 'Class 5, Class 6, Class 7'
$ws2.Range('B13').Value = 'iOS'
$ws2.Range('C13').Value = 'verifyPrimeClasses_Parent'
$ws2.Range('D13').Value = 'Parent'
$ws2.Range('E13').Value = 'Prime Classes'
$ws2.Range('F13').Value = 'Class 5, Class 6, Class 7'
$ws2.Range('B14').Value = 'Web'
$ws2.Range('C14').Value = 'verifyPrimeSubjects_Old'
$ws2.Range('D14').Value = 'Admin'
$ws2.Range('E14').Value = 'Prime Subjects'
$ws2.Range('F14').Value = 'English, Hindi'
$ws2.Range('B15').Value = 'Android'
$ws2.Range('C15').Value = 'verifyPrimeSubjects_Old'
$ws2.Range('D15').Value = 'Admin'
$ws2.Range('E15').Value = 'Prime Subjects'
$ws2.Range('F15').Value = 'English, Hindi'
$ws2.Range('B16').Value = 'iOS'
$ws2.Range('C16').Value = 'verifyPrimeSubjects_Old'
$ws2.Range('D16').Value = 'Admin'
$ws2.Range('E16').Value = 'Prime Subjects'
$ws2.Range('F16').Value = 'English, Hindi'
$ws2.Range('B17').Value = 'Web'
$ws2.Range('C17').Value = 'verifyPrimeSubjects_Old'
$ws2.Range('D17').Value = 'Principal'
$ws2.Range('E17').Value = 'Prime Subjects'
$ws2.Range('F17').Value = 'English, Hindi'
$ws2.Range('B18').Value = 'Android'
$ws2.Range('C18').Value = 'verifyPrimeSubjects_Old'
$ws2.Range('D18').Value = 'Principal'
$ws2.Range('E18').Value = 'Prime Subjects'
$ws2.Range('F18').Value = 'English, Hindi'
$ws2.Range('B19').Value = 'iOS'
$ws2.Range('C19').Value = 'verifyPrimeSubjects_Old'
$ws2.Range('D19').Value = 'Principal'
$ws2.Range('E19').Value = 'Prime Subjects'
$ws2.Range('F19').Value = 'English, Hindi'
$ws2.Range('B20').Value = 'Web'
$ws2.Range('C20').Value = 'verifyPrimeSubjects_Old'
$ws2.Range('D20').Value = 'Teacher'
$ws2.Range('E20').Value = 'Prime Subjects'
$ws2.Range('F20').Value = 'English, Hindi'
$ws2.Range('B21').Value = 'Android'
$ws2.Range('C21').Value = 'verifyPrimeSubjects_Old'
$ws2.Range('D21').Value = 'Teacher'
$ws2.Range('E21').Value = 'Prime Subjects'
$ws2.Range('F21').Value = 'English, Hindi'
$ws2.Range('B22').Value = 'iOS'
$ws2.Range('C22').Value = 'verifyPrimeSubjects_Old'
$ws2.Range('D22').Value = 'Teacher'
$ws2.Range('E22').Value = 'Prime Subjects'
$ws2.Range('F22').Value = 'English, Hindi'
$ws2.Range('B23').Value = 'Web'
$ws2.Range('C23').Value = 'verifyPrimeSubjects_Old'
$ws2.Range('D23').Value = 'Parent'
$ws2.Range('E23').Value = 'Prime Subjects'
$ws2.Range('F23').Value = 'English, Hindi'
$ws2.Range('B24').Value = 'Android'
$ws2.Range('C24').Value = 'verifyPrimeSubjects_Old'
$ws2.Range('D24').Value = 'Parent'
$ws2.Range('E24').Value = 'Prime Subjects'
$ws2.Range('F24').Value = 'English, Hindi'
$ws2.Range('B25').Value = 'iOS'
$ws2.Range('C25').Value = 'verifyPrimeSubjects_Old'
$ws2.Range('D25').Value = 'Parent'
$ws2.Range('E25').Value = 'Prime Subjects'
$ws2.Range('F25').Value = 'English, Hindi'
$ws2.Range('B26').Value = 'Web'
$ws2.Range('C26').Value = 'verifyPrimeSubjects_Old'
$ws2.Range('D26').Value = 'Student'
$ws2.Range('E26').Value = 'Prime Subjects'
$ws2.Range('F26').Value = 'English, Hindi'
$ws2.Range('B27').Value = 'Android'
$ws2.Range('C27').Value = 'verifyPrimeSubjects_Old'
$ws2.Range('D27').Value = 'Student'
$ws2.Range('E27').Value = 'Prime Subjects'
$ws2.Range('F27').Value = 'English, Hindi'
$ws2.Range('B28').Value = 'iOS'
$ws2.Range('C28').Value = 'verifyPrimeSubjects_Old'
$ws2.Range('D28').Value = 'Student'
$ws2.Range('E28').Value = 'Prime Subjects'
$ws2.Range('F28').Value = 'English, Hindi'
$ws2.Range('B29').Value = 'Web'
$ws2.Range('C29').Value = 'verifyPrimeSubjects_Old'
$ws2.Range('D29').Value = 'Guest'
$ws2.Range('E29').Value = 'Prime Subjects'
$ws2.Range('F29').Value = 'English, Hindi'
$ws2.Range('B30').Value = 'Web'
$ws2.Range('C30').Value = 'verifyPrimeSubjects_Old'
$ws2.Range('D30').Value = 'Guest'
$ws2.Range('E30').Value = 'Prime Subjects'
$ws2.Range('F30').Value = 'English, Hindi'
$ws2.Range('B31').Value = 'Android'
$ws2.Range('C31').Value = 'verifyPrimeSubjects_Old'
$ws2.Range('D31').Value = 'Guest'
$ws2.Range('E31').Value = 'Prime Subjects'
$ws2.Range('F31').Value = 'English, Hindi'
$ws2.Range('B32').Value = 'Web'
$ws2.Range('C32').Value = 'verifyQuizDashboardClasses'
$ws2.Range('D32').Value = 'Admin'
$ws2.Range('E32').Value = 'Quiz Dashboard Classes'
$ws2.Range('F32').Value = 'Class 3-A, Class 3-A, Class 3-A, Class 3-A, Class 3-B, Class 3-B, Class 3-B, Class 3-B, Class 3-C, Class 3-C, Class 3-C, Class 4-A, Class 4-A, Class 4-A, Class 4-A, Class 4-B, Class 4-B, Class 4-B, Class 4-B, Class 4-C, Class 4-C, Class 4-C, Class 4-C, Class 5-A, Class 5-A, Class 5-A, Class 5-A, Class 5-B, Class 5-B, Class 5-B, Class 5-B, Class 5-C, Class 5-C, Class 5-C, Class 5-C, Class 6-A, Class 6-A, Class 6-A, Class 6-A, Class 6-B, Class 6-B, Class 6-B, Class 6-B, Class 6-C, Class 6-C, Class 6-C, Class 6-C, Class 7-A, Class 7-A, Class 7-A, Class 7-A, Class 7-B, Class 7-B, Class 7-B, Class 7-B, Class 7-C, Class 7-C, Class 7-C, Class 7-C, Class 8-A, Class 8-A, Class 8-A, Class 8-A, Class 8-B, Class 8-B, Class 8-B, Class 8-B, Class 8-C, Class 8-C, Class 8-C, Class 8-C, Class 9-A, Class 9-A, Class 9-A, Class 9-A, Class 9-B, Class 9-B, Class 9-B, Class 9-B, Class 9-C, Class 9-C, Class 9-C, Class 9-C, Class 10-A, Class 10-A, Class 10-B, Class 10-B, Class 10-C, Class 10-C, Class 11-A, Class 11-A, Class 11-A, Class 11-A, Class 11-B, Class 11-B, Class 11-B, Class 11-B, Class 11-C, Class 11-C, Class 11-C, Class 11-C, Class 12-A, Class 12-A, Class 12-A, Class 12-A, Class 12-B, Class 12-B, Class 12-B, Class 12-B, Class 12-C, Class 12-C, Class 12-C'
$ws2.Range('B33').Value = 'Android'
$ws2.Range('C33').Value = 'verifyQuizDashboardClasses'
$ws2.Range('D33').Value = 'Admin'
$ws2.Range('E33').Value = 'Quiz Dashboard Classes'
$ws2.Range('F33').Value = 'Class 3-A, Class 3-A, Class 3-A, Class 3-A, Class 3-B, Class 3-B, Class 3-B, Class 3-B, Class 3-C, Class 3-C, Class 3-C, Class 4-A, Class 4-A, Class 4-A, Class 4-A, Class 4-B, Class 4-B, Class 4-B, Class 4-B, Class 4-C, Class 4-C, Class 4-C, Class 4-C, Class 5-A, Class 5-A, Class 5-A, Class 5-A, Class 5-B, Class 5-B, Class 5-B, Class 5-B, Class 5-C, Class 5-C, Class 5-C, Class 5-C, Class 6-A, Class 6-A, Class 6-A, Class 6-A, Class 6-B, Class 6-B, Class 6-B, Class 6-B, Class 6-C, Class 6-C, Class 6-C, Class 6-C, Class 7-A, Class 7-A, Class 7-A, Class 7-A, Class 7-B, Class 7-B, Class 7-B, Class 7-B, Class 7-C, Class 7-C, Class 7-C, Class 7-C, Class 8-A, Class 8-A, Class 8-A, Class 8-A, Class 8-B, Class 8-B, Class 8-B, Class 8-B, Class 8-C, Class 8-C, Class 8-C, Class 8-C, Class 9-A, Class 9-A, Class 9-A, Class 9-A, Class 9-B, Class 9-B, Class 9-B, Class 9-B, Class 9-C, Class 9-C, Class 9-C, Class 9-C, Class 10-A, Class 10-A, Class 10-B, Class 10-B, Class 10-C, Class 10-C, Class 11-A, Class 11-A, Class 11-A, Class 11-A, Class 11-B, Class 11-B, Class 11-B, Class 11-B, Class 11-C, Class 11-C, Class 11-C, Class 11-C, Class 12-A, Class 12-A, Class 12-A, Class 12-A, Class 12-B, Class 12-B, Class 12-B, Class 12-B, Class 12-C, Class 12-C, Class 12-C'
$ws2.Range('B34').Value = 'iOS'
$ws2.Range('C34').Value = 'verifyQuizDashboardClasses'
$ws2.Range('D34').Value = 'Admin'
$ws2.Range('E34').Value = 'Quiz Dashboard Classes'
$ws2.Range('F34').Value = 'Class 3-A, Class 3-A, Class 3-A, Class 3-A, Class 3-B, Class 3-B, Class 3-B, Class 3-B, Class 3-C, Class 3-C, Class 3-C, Class 4-A, Class 4-A, Class 4-A, Class 4-A, Class 4-B, Class 4-B, Class 4-B, Class 4-B, Class 4-C, Class 4-C, Class 4-C, Class 4-C, Class 5-A, Class 5-A, Class 5-A, Class 5-A, Class 5-B, Class 5-B, Class 5-B, Class 5-B, Class 5-C, Class 5-C, Class 5-C, Class 5-C, Class 6-A, Class 6-A, Class 6-A, Class 6-A, Class 6-B, Class 6-B, Class 6-B, Class 6-B, Class 6-C, Class 6-C, Class 6-C, Class 6-C, Class 7-A, Class 7-A, Class 7-A, Class 7-A, Class 7-B, Class 7-B, Class 7-B, Class 7-B, Class 7-C, Class 7-C, Class 7-C, Class 7-C, Class 8-A, Class 8-A, Class 8-A, Class 8-A, Class 8-B, Class 8-B, Class 8-B, Class 8-B, Class 8-C, Class 8-C, Class 8-C, Class 8-C, Class 9-A, Class 9-A, Class 9-A, Class 9-A, Class 9-B, Class 9-B, Class 9-B, Class 9-B, Class 9-C, Class 9-C, Class 9-C, Class 9-C, Class 10-A, Class 10-A, Class 10-B, Class 10-B, Class 10-C, Class 10-C, Class 11-A, Class 11-A, Class 11-A, Class 11-A, Class 11-B, Class 11-B, Class 11-B, Class 11-B, Class 11-C, Class 11-C, Class 11-C, Class 11-C, Class 12-A, Class 12-A, Class 12-A, Class 12-A, Class 12-B, Class 12-B, Class 12-B, Class 12-B, Class 12-C, Class 12-C, Class 12-C'
$ws2.Range('B35').Value = 'Web'
$ws2.Range('C35').Value = 'verifyQuizDashboardClasses'
$ws2.Range('D35').Value = 'Principal'
$ws2.Range('E35').Value = 'Quiz Dashboard Classes'
$ws2.Range('F35').Value = 'Class 3-A, Class 3-A, Class 3-A, Class 3-A, Class 3-B, Class 3-B, Class 3-B, Class 3-B, Class 3-C, Class 3-C, Class 3-C, Class 4-A, Class 4-A, Class 4-A, Class 4-A, Class 4-B, Class 4-B, Class 4-B, Class 4-B, Class 4-C, Class 4-C, Class 4-C, Class 4-C, Class 5-A, Class 5-A, Class 5-A, Class 5-A, Class 5-B, Class 5-B, Class 5-B, Class 5-B, Class 5-C, Class 5-C, Class 5-C, Class 5-C, Class 6-A, Class 6-A, Class 6-A, Class 6-A, Class 6-B, Class 6-B, Class 6-B, Class 6-B, Class 6-C, Class 6-C, Class 6-C, Class 6-C, Class 7-A, Class 7-A, Class 7-A, Class 7-A, Class 7-B, Class 7-B, Class 7-B, Class 7-B, Class 7-C, Class 7-C, Class 7-C, Class 7-C, Class 8-A, Class 8-A, Class 8-A, Class 8-A, Class 8-B, Class 8-B, Class 8-B, Class 8-B, Class 8-C, Class 8-C, Class 8-C, Class 8-C, Class 9-A, Class 9-A, Class 9-A, Class 9-A, Class 9-B, Class 9-B, Class 9-B, Class 9-B, Class 9-C, Class 9-C, Class 9-C, Class 9-C, Class 10-A, Class 10-A, Class 10-B, Class 10-B, Class 10-C, Class 10-C, Class 11-A, Class 11-A, Class 11-A, Class 11-A, Class 11-B, Class 11-B, Class 11-B, Class 11-B, Class 11-C, Class 11-C, Class 11-C, Class 11-C, Class 12-A, Class 12-A, Class 12-A, Class 12-A, Class 12-B, Class 12-B, Class 12-B, Class 12-B, Class 12-C, Class 12-C, Class 12-C'
$ws2.Range('B36').Value = 'Android'
$ws2.Range('C36').Value = 'verifyQuizDashboardClasses'
$ws2.Range('D36').Value = 'Principal'
$ws2.Range('E36').Value = 'Quiz Dashboard Classes'
$ws2.Range('F36').Value = 'Class 3-A, Class 3-A, Class 3-A, Class 3-A, Class 3-B, Class 3-B, Class 3-B, Class 3-B, Class 3-C, Class 3-C, Class 3-C, Class 4-A, Class 4-A, Class 4-A, Class 4-A, Class 4-B, Class 4-B, Class 4-B, Class 4-B, Class 4-C, Class 4-C, Class 4-C, Class 4-C, Class 5-A, Class 5-A, Class 5-A, Class 5-A, Class 5-B, Class 5-B, Class 5-B, Class 5-B, Class 5-C, Class 5-C, Class 5-C, Class 5-C, Class 6-A, Class 6-A, Class 6-A, Class 6-A, Class 6-B, Class 6-B, Class 6-B, Class 6-B, Class 6-C, Class 6-C, Class 6-C, Class 6-C, Class 7-A, Class 7-A, Class 7-A, Class 7-A, Class 7-B, Class 7-B, Class 7-B, Class 7-B, Class 7-C, Class 7-C, Class 7-C, Class 7-C, Class 8-A, Class 8-A, Class 8-A, Class 8-A, Class 8-B, Class 8-B, Class 8-B, Class 8-B, Class 8-C, Class 8-C, Class 8-C, Class 8-C, Class 9-A, Class 9-A, Class 9-A, Class 9-A, Class 9-B, Class 9-B, Class 9-B, Class 9-B, Class 9-C, Class 9-C, Class 9-C, Class 9-C, Class 10-A, Class 10-A, Class 10-B, Class 10-B, Class 10-C, Class 10-C, Class 11-A, Class 11-A, Class 11-A, Class 11-A, Class 11-B, Class 11-B, Class 11-B, Class 11-B, Class 11-C, Class 11-C, Class 11-C, Class 11-C, Class 12-A, Class 12-A, Class 12-A, Class 12-A, Class 12-B, Class 12-B, Class 12-B, Class 12-B, Class 12-C, Class 12-C, Class 12-C'
$ws2.Range('B37').Value = 'iOS'
$ws2.Range('C37').Value = 'verifyQuizDashboardClasses'
$ws2.Range('D37').Value = 'Principal'
$ws2.Range('E37').Value = 'Quiz Dashboard Classes'
$ws2.Range('F37').Value = 'Class 3-A, Class 3-A, Class 3-A, Class 3-A, Class 3-B, Class 3-B, Class 3-B, Class 3-B, Class 3-C, Class 3-C, Class 3-C, Class 4-A, Class 4-A, Class 4-A, Class 4-A, Class 4-B, Class 4-B, Class 4-B, Class 4-B, Class 4-C, Class 4-C, Class 4-C, Class 4-C, Class 5-A, Class 5-A, Class 5-A, Class 5-A, Class 5-B, Class 5-B, Class 5-B, Class 5-B, Class 5-C, Class 5-C, Class 5-C, Class 5-C, Class 6-A, Class 6-A, Class 6-A, Class 6-A, Class 6-B, Class 6-B, Class 6-B, Class 6-B, Class 6-C, Class 6-C, Class 6-C, Class 6-C, Class 7-A, Class 7-A, Class 7-A, Class 7-A, Class 7-B, Class 7-B, Class 7-B, Class 7-B, Class 7-C, Class 7-C, Class 7-C, Class 7-C, Class 8-A, Class 8-A, Class 8-A, Class 8-A, Class 8-B, Class 8-B, Class 8-B, Class 8-B, Class 8-C, Class 8-C, Class 8-C, Class 8-C, Class 9-A, Class 9-A, Class 9-A, Class 9-A, Class 9-B, Class 9-B, Class 9-B, Class 9-B, Class 9-C, Class 9-C, Class 9-C, Class 9-C, Class 10-A, Class 10-A, Class 10-B, Class 10-B, Class 10-C, Class 10-C, Class 11-A, Class 11-A, Class 11-A, Class 11-A, Class 11-B, Class 11-B, Class 11-B, Class 11-B, Class 11-C, Class 11-C, Class 11-C, Class 11-C, Class 12-A, Class 12-A, Class 12-A, Class 12-A, Class 12-B, Class 12-B, Class 12-B, Class 12-B, Class 12-C, Class 12-C, Class 12-C'
$ws2.Range('B38').Value = 'Web'
$ws2.Range('C38').Value = 'verifyQuizDashboardClasses'
$ws2.Range('D38').Value = 'Teacher'
$ws2.Range('E38').Value = 'Quiz Dashboard Classes'
$ws2.Range('F38').Value = 'Class 5-C, Class 6-A'
$ws2.Range('B39').Value = 'Android'
$ws2.Range('C39').Value = 'verifyQuizDashboardClasses'
$ws2.Range('D39').Value = 'Teacher'
$ws2.Range('E39').Value = 'Quiz Dashboard Classes'
$ws2.Range('F39').Value = 'Class 5-C, Class 6-A'
$ws2.Range('B40').Value = 'iOS'
$ws2.Range('C40').Value = 'verifyQuizDashboardClasses'
$ws2.Range('D40').Value = 'Teacher'
$ws2.Range('E40').Value = 'Quiz Dashboard Classes'
$ws2.Range('F40').Value = 'Class 5-C, Class 6-A'
$ws2.Range('B41').Value = 'Web'
$ws2.Range('C41').Value = 'searchAndViewContentSchool_Old'
$ws2.Range('D41').Value = 'Admin'
$ws2.Range('E41').Value = 'Search String'
$ws2.Range('F41').Value = 'Autotrophic Nutrition'
$ws2.Range('B42').Value = 'Android'
$ws2.Range('C42').Value = 'searchAndViewContentSchool_Old'
$ws2.Range('D42').Value = 'Admin'
$ws2.Range('E42').Value = 'Search String'
$ws2.Range('F42').Value = 'Autotrophic Nutrition'
$ws2.Range('B43').Value = 'iOS'
$ws2.Range('C43').Value = 'searchAndViewContentSchool_Old'
$ws2.Range('D43').Value = 'Admin'
$ws2.Range('E43').Value = 'Search String'
$ws2.Range('F43').Value = 'Autotrophic Nutrition'
$ws2.Range('B44').Value = 'Web'
$ws2.Range('C44').Value = 'searchAndViewContentSchool_Old'
$ws2.Range('D44').Value = 'Principal'
$ws2.Range('E44').Value = 'Search String'
$ws2.Range('F44').Value = 'Autotrophic Nutrition'
$ws2.Range('B45').Value = 'Android'
$ws2.Range('C45').Value = 'searchAndViewContentSchool_Old'
$ws2.Range('D45').Value = 'Principal'
$ws2.Range('E45').Value = 'Search String'
$ws2.Range('F45').Value = 'Autotrophic Nutrition'
$ws2.Range('B46').Value = 'iOS'
$ws2.Range('C46').Value = 'searchAndViewContentSchool_Old'
$ws2.Range('D46').Value = 'Principal'
$ws2.Range('E46').Value = 'Search String'
$ws2.Range('F46').Value = 'Autotrophic Nutrition'
$ws2.Range('B47').Value = 'Web'
$ws2.Range('C47').Value = 'searchAndViewContentSchool_Old'
$ws2.Range('D47').Value = 'Teacher'
$ws2.Range('E47').Value = 'Search String'
$ws2.Range('F47').Value = 'Autotrophic Nutrition'
$ws2.Range('B48').Value = 'Android'
$ws2.Range('C48').Value = 'searchAndViewContentSchool_Old'
$ws2.Range('D48').Value = 'Teacher'
$ws2.Range('E48').Value = 'Search String'
$ws2.Range('F48').Value = 'Autotrophic Nutrition'
$ws2.Range('B49').Value = 'iOS'
$ws2.Range('C49').Value = 'searchAndViewContentSchool_Old'
$ws2.Range('D49').Value = 'Teacher'
$ws2.Range('E49').Value = 'Search String'
$ws2.Range('F49').Value = 'Autotrophic Nutrition'
$ws2.Range('B50').Value = 'Web'
$ws2.Range('C50').Value = 'searchAndViewContentStudent_Old'
$ws2.Range('D50').Value = 'Parent'
$ws2.Range('E50').Value = 'Search String'
$ws2.Range('F50').Value = 'Making stone tools and the discovery of fire'
$ws2.Range('B51').Value = 'Android'
$ws2.Range('C51').Value = 'searchAndViewContentStudent_Old'
$ws2.Range('D51').Value = 'Parent'
$ws2.Range('E51').Value = 'Search String'
$ws2.Range('F51').Value = 'Making stone tools and the discovery of fire'
$ws2.Range('B52').Value = 'iOS'
$ws2.Range('C52').Value = 'searchAndViewContentStudent_Old'
$ws2.Range('D52').Value = 'Parent'
$ws2.Range('E52').Value = 'Search String'
$ws2.Range('F52').Value = 'Making stone tools and the discovery of fire'
$ws2.Range('B53').Value = 'Web'
$ws2.Range('C53').Value = 'searchAndViewContentStudent_Old'
$ws2.Range('D53').Value = 'Student'
$ws2.Range('E53').Value = 'Search String'
$ws2.Range('F53').Value = 'Making stone tools and the discovery of fire'
$ws2.Range('B54').Value = 'Android'
$ws2.Range('C54').Value = 'searchAndViewContentStudent_Old'
$ws2.Range('D54').Value = 'Student'
$ws2.Range('E54').Value = 'Search String'
$ws2.Range('F54').Value = 'Making stone tools and the discovery of fire'
$ws2.Range('B55').Value = 'iOS'
$ws2.Range('C55').Value = 'searchAndViewContentStudent_Old'
$ws2.Range('D55').Value = 'Student'
$ws2.Range('E55').Value = 'Search String'
$ws2.Range('F55').Value = 'Making stone tools and the discovery of fire'
$ws2.Range('B56').Value = 'Web'
$ws2.Range('C56').Value = 'searchAndViewContentStudent_Old'
$ws2.Range('D56').Value = 'Guest'
$ws2.Range('E56').Value = 'Search String'
$ws2.Range('F56').Value = 'Making stone tools and the discovery of fire'
$ws2.Range('B57').Value = 'Android'
$ws2.Range('C57').Value = 'searchAndViewContentStudent_Old'
$ws2.Range('D57').Value = 'Guest'
$ws2.Range('E57').Value = 'Search String'
$ws2.Range('F57').Value = 'Making stone tools and the discovery of fire'
$ws2.Range('B58').Value = 'iOS'
$ws2.Range('C58').Value = 'searchAndViewContentStudent_Old'
$ws2.Range('D58').Value = 'Guest'
$ws2.Range('E58').Value = 'Search String'
$ws2.Range('F58').Value = 'Making stone tools and the discovery of fire'
$ws2.Range('B59').Value = 'Web'
$ws2.Range('C59').Value = 'verifyCountOfDoubtTabs'
$ws2.Range('D59').Value = 'Teacher'
$ws2.Range('E59').Value = 'DoubtForum Tabs'
$ws2.Range('F59').Value = 'Doubts, My Doubts, Posts For Me'
$ws2.Range('B60').Value = 'Android'
$ws2.Range('C60').Value = 'verifyCountOfDoubtTabs'
$ws2.Range('D60').Value = 'Teacher'
$ws2.Range('E60').Value = 'DoubtForum Tabs'
$ws2.Range('F60').Value = 'Doubts, My Doubts, Post For Me'
$ws2.Range('B61').Value = 'Web'
$ws2.Range('C61').Value = 'verifyCountOfDoubtTabs'
$ws2.Range('D61').Value = 'Student'
$ws2.Range('E61').Value = 'DoubtForum Tabs'
$ws2.Range('F61').Value = 'Doubts, My Doubts'
$ws2.Range('B62').Value = 'Web'
$ws2.Range('C62').Value = 'verifyCountOfDoubtTabs'
$ws2.Range('D62').Value = 'Guest'
$ws2.Range('E62').Value = 'DoubtForum Tabs'
$ws2.Range('F62').Value = 'Doubts, My Doubts'

# Apply explicit (black) font color to E62 on the staging sheet, matching
# the s="1" style used in the source workbook for that cell.
$ws2.Range('E62').Font.Color = 0

# Restore view/selection state (per-sheet selection, mirroring the
# cursor positions recorded in the target workbook). "production" is the
# active tab, so select its target cell last.
$ws2.Range("C53").Select()
$ws1.Range("C69").Select()
